$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows after the existing "GET /redfish/v1/AccountService" row
# (shifts old rows 3-12 down to 5-14), then populate them with new
# GET requests used to look up a user's id by name.
$ws.Rows("3:4").Insert()

$ws.Range("A3").Value = "GET"
$ws.Range("B3").Value = "/redfish/v1/AccountService/Accounts/1"

$ws.Range("A4").Value = "GET"
$ws.Range("B4").Value = "/redfish/v1/AccountService/Accounts/admin"

# The DELAY row (now row 14) changes its wait value from 10 to 5.
$ws.Range("B14").Value = 5

# Leave the cursor where the author ended up after editing.
$ws.Range("C34").Select() | Out-Null
